$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.038012492658291
$ws.Range("D2").Value = 1.046049527104035
$ws.Range("E2").Value = 1.046453301742122
$ws.Range("F2").Value = 1.056661513366447
$ws.Range("I2").Value = 1.038012640255427
$ws.Range("J2").Value = 1.043112241797363
$ws.Range("K2").Value = 1.048816052034104
$ws.Range("L2").Value = 1.049218695171206
$ws.Range("M2").Value = 1.059398611523202
$ws.Range("N2").Value = 1.018316225590214

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.038865825066547
$ws.Range("D3").Value = 1.046800092458259
$ws.Range("E3").Value = 1.047203905475387
$ws.Range("F3").Value = 1.057476800197066
$ws.Range("I3").Value = 1.038135541311707
$ws.Range("J3").Value = 1.043610673772147
$ws.Range("K3").Value = 1.049378340121172
$ws.Range("L3").Value = 1.049781103153058
$ws.Range("M3").Value = 1.060027577364837
$ws.Range("N3").Value = 1.018481233640325

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.039418332098627
$ws.Range("D4").Value = 1.047286369663836
$ws.Range("E4").Value = 1.047690270150612
$ws.Range("F4").Value = 1.058005005064573
$ws.Range("I4").Value = 1.038213205078084
$ws.Range("J4").Value = 1.043932875261951
$ws.Range("K4").Value = 1.049742123408819
$ws.Range("L4").Value = 1.050145026446113
$ws.Range("M4").Value = 1.06043456611207
$ws.Range("N4").Value = 1.01858788565255

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.039650686737884
$ws.Range("D5").Value = 1.047490945059368
$ws.Range("E5").Value = 1.047894897312929
$ws.Range("F5").Value = 1.05822721835735
$ws.Range("I5").Value = 1.038245408365186
$ws.Range("J5").Value = 1.044068251427904
$ws.Range("K5").Value = 1.049895043304334
$ws.Range("L5").Value = 1.050298020149936
$ws.Range("M5").Value = 1.060605663892511
$ws.Range("N5").Value = 1.018632693136734

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.039689704787554
$ws.Range("D6").Value = 1.047525302585013
$ws.Range("E6").Value = 1.047929264409716
$ws.Range("F6").Value = 1.058264538035137
$ws.Range("I6").Value = 1.038250789225171
$ws.Range("J6").Value = 1.044090977106602
$ws.Range("K6").Value = 1.04992071832782
$ws.Range("L6").Value = 1.050323708442992
$ws.Range("M6").Value = 1.06063439189696
$ws.Range("N6").Value = 1.018640214794713

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.039421436517732
$ws.Range("D7").Value = 1.047289102645025
$ws.Range("E7").Value = 1.047693003764594
$ws.Range("F7").Value = 1.05800797367913
$ws.Range("I7").Value = 1.038213637135834
$ws.Range("J7").Value = 1.043934684470533
$ws.Range("K7").Value = 1.049744166791206
$ws.Range("L7").Value = 1.050147070755981
$ws.Range("M7").Value = 1.060436852332237
$ws.Range("N7").Value = 1.018588484486986

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.038300808143358
$ws.Range("D8").Value = 1.046303056496175
$ws.Range("E8").Value = 1.046706831135089
$ws.Range("F8").Value = 1.056936905566863
$ws.Range("I8").Value = 1.038054560108408
$ws.Range("J8").Value = 1.043280754098379
$ws.Range("K8").Value = 1.049006090350796
$ws.Range("L8").Value = 1.049408761084429
$ws.Range("M8").Value = 1.059611171269992
$ws.Range("N8").Value = 1.018372015233897

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.036328824309272
$ws.Range("D9").Value = 1.044570280998535
$ws.Range("E9").Value = 1.044974312513483
$ws.Range("F9").Value = 1.055054687203111
$ws.Range("I9").Value = 1.037760030461669
$ws.Range("J9").Value = 1.042126075265809
$ws.Range("K9").Value = 1.047705152557897
$ws.Range("L9").Value = 1.048107890188721
$ws.Range("M9").Value = 1.058156333856942
$ws.Range("N9").Value = 1.01798967662952

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.035016085327877
$ws.Range("D10").Value = 1.043418411103697
$ws.Range("E10").Value = 1.043822936505795
$ws.Range("F10").Value = 1.053803448613586
$ws.Range("I10").Value = 1.037554174104555
$ws.Range("J10").Value = 1.04135477690378
$ws.Range("K10").Value = 1.046837712284956
$ws.Range("L10").Value = 1.047240815536005
$ws.Range("M10").Value = 1.057186620684846
$ws.Range("N10").Value = 1.017734212636666

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.034448129778071
$ws.Range("D11").Value = 1.042920447991769
$ws.Range("E11").Value = 1.043325263804723
$ws.Range("F11").Value = 1.053262520778113
$ws.Range("I11").Value = 1.037462793322317
$ws.Range("J11").Value = 1.041020454397112
$ws.Range("K11").Value = 1.046462084475989
$ws.Range("L11").Value = 1.046865422136699
$ws.Range("M11").Value = 1.05676678725681
$ws.Range("N11").Value = 1.017623464244392

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.034237237880507
$ws.Range("D12").Value = 1.042735605009627
$ws.Range("E12").Value = 1.043140540199601
$ws.Range("F12").Value = 1.053061728295397
$ws.Range("I12").Value = 1.037428514231593
$ws.Range("J12").Value = 1.040896221613709
$ws.Range("K12").Value = 1.046322557997943
$ws.Range("L12").Value = 1.046725994158627
$ws.Range("M12").Value = 1.056610852919235
$ws.Range("N12").Value = 1.017582308180511

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.0342824716208
$ws.Range("D13").Value = 1.042775248874773
$ws.Range("E13").Value = 1.043180157937257
$ws.Range("F13").Value = 1.053104792950422
$ws.Range("I13").Value = 1.037435882412584
$ws.Range("J13").Value = 1.040922872231956
$ws.Range("K13").Value = 1.046352486949019
$ws.Range("L13").Value = 1.04675590146394
$ws.Range("M13").Value = 1.056644300862648
$ws.Range("N13").Value = 1.017591137156309

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.034430695908301
$ws.Range("D14").Value = 1.042905166299558
$ws.Range("E14").Value = 1.043309991745522
$ws.Range("F14").Value = 1.053245920500562
$ws.Range("I14").Value = 1.037459966655284
$ws.Range("J14").Value = 1.041010186299271
$ws.Range("K14").Value = 1.04645055119754
$ws.Range("L14").Value = 1.046853896767139
$ws.Range("M14").Value = 1.056753897452233
$ws.Range("N14").Value = 1.01762006265986

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.034522031383404
$ws.Range("D15").Value = 1.042985229026286
$ws.Range("E15").Value = 1.043390004478063
$ws.Range("F15").Value = 1.053332891419035
$ws.Range("I15").Value = 1.037474761228616
$ws.Range("J15").Value = 1.04106397672388
$ws.Range("K15").Value = 1.046510971651947
$ws.Range("L15").Value = 1.046914276256542
$ws.Range("M15").Value = 1.056821424967906
$ws.Range("N15").Value = 1.017637882089791

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.035053788687564
$ws.Range("D16").Value = 1.043451476378166
$ws.Range("E16").Value = 1.043855984119913
$ws.Range("F16").Value = 1.053839366643009
$ws.Range("I16").Value = 1.037560191550758
$ws.Range("J16").Value = 1.041376957629101
$ws.Range("K16").Value = 1.046862641161785
$ws.Range("L16").Value = 1.047265730458836
$ws.Range("M16").Value = 1.057214485026953
$ws.Range("N16").Value = 1.017741559927203

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.035387472719255
$ws.Range("D17").Value = 1.043744157643849
$ws.Range("E17").Value = 1.044148517956821
$ws.Range("F17").Value = 1.054157298620831
$ws.Range("I17").Value = 1.037613179675932
$ws.Range("J17").Value = 1.041573190808064
$ws.Range("K17").Value = 1.047083229731447
$ws.Range("L17").Value = 1.047486204342475
$ws.Range("M17").Value = 1.057461058127077
$ws.Range("N17").Value = 1.017806559610258

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.03558215012243
$ws.Range("D18").Value = 1.04391495113816
$ws.Range("E18").Value = 1.044319232844543
$ws.Range("F18").Value = 1.05434282645877
$ws.Range("I18").Value = 1.037643870217261
$ws.Range("J18").Value = 1.041687616901928
$ws.Range("K18").Value = 1.047211893239842
$ws.Range("L18").Value = 1.047614808298551
$ws.Range("M18").Value = 1.057604885604553
$ws.Range("N18").Value = 1.017844460174123

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.035648537695832
$ws.Range("D19").Value = 1.043973200388109
$ws.Range("E19").Value = 1.044377456544369
$ws.Range("F19").Value = 1.05440610075867
$ws.Range("I19").Value = 1.037654298143723
$ws.Range("J19").Value = 1.041726627536709
$ws.Range("K19").Value = 1.047255763774444
$ws.Range("L19").Value = 1.047658659772773
$ws.Range("M19").Value = 1.057653927949985
$ws.Range("N19").Value = 1.017857381126686

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.035351666898604
$ws.Range("D20").Value = 1.043712747687955
$ws.Range("E20").Value = 1.04411712305463
$ws.Range("F20").Value = 1.054123178877947
$ws.Range("I20").Value = 1.037607516942499
$ws.Range("J20").Value = 1.041552140290004
$ws.Range("K20").Value = 1.047059562869139
$ws.Range("L20").Value = 1.047462549025483
$ws.Range("M20").Value = 1.057434602579888
$ws.Range("N20").Value = 1.01779958706462

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.034387045547323
$ws.Range("D21").Value = 1.042866905454212
$ws.Range("E21").Value = 1.043271755205517
$ws.Range("F21").Value = 1.053204358289283
$ws.Range("I21").Value = 1.037452883719806
$ws.Range("J21").Value = 1.040984475866456
$ws.Range("K21").Value = 1.046421673752373
$ws.Range("L21").Value = 1.046825039309152
$ws.Range("M21").Value = 1.056721623685321
$ws.Range("N21").Value = 1.017611545345699

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.033780967655301
$ws.Range("D22").Value = 1.04233580163724
$ws.Range("E22").Value = 1.04274101622227
$ws.Range("F22").Value = 1.052627425164462
$ws.Range("I22").Value = 1.037353714960701
$ws.Range("J22").Value = 1.040627271090938
$ws.Range("K22").Value = 1.046020599109267
$ws.Range("L22").Value = 1.046424269338211
$ws.Range("M22").Value = 1.056273407070078
$ws.Range("N22").Value = 1.017493205316469

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.03410222087124
$ws.Range("D23").Value = 1.042617281822636
$ws.Range("E23").Value = 1.043022296691047
$ws.Range("F23").Value = 1.052933195075597
$ws.Range("I23").Value = 1.037406470240368
$ws.Range("J23").Value = 1.040816659294575
$ws.Range("K23").Value = 1.046233216670754
$ws.Range("L23").Value = 1.046636719119844
$ws.Range("M23").Value = 1.056511008794557
$ws.Range("N23").Value = 1.017555949944831

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.035367845868368
$ws.Range("D24").Value = 1.043726940255872
$ws.Range("E24").Value = 1.044131308797606
$ws.Range("F24").Value = 1.054138595863341
$ws.Range("I24").Value = 1.037610076357137
$ws.Range("J24").Value = 1.041561652216357
$ws.Range("K24").Value = 1.047070256912062
$ws.Range("L24").Value = 1.047473237828889
$ws.Range("M24").Value = 1.05744655668591
$ws.Range("N24").Value = 1.017802737697009

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.03683829863701
$ws.Range("D25").Value = 1.04501766901809
$ws.Range("E25").Value = 1.04542157726353
$ws.Range("F25").Value = 1.055540664463545
$ws.Range("I25").Value = 1.037837852421024
$ws.Range("J25").Value = 1.042424859914474
$ws.Range("K25").Value = 1.048041508800775
$ws.Range("L25").Value = 1.048444172534656
$ws.Range("M25").Value = 1.058532419978219
$ws.Range("N25").Value = 1.018088623129696
